$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.782084584236145
$ws.Range("B1").Value = 3.173549652099609
$ws.Range("C1").Value = 3.627765893936157
$ws.Range("D1").Value = 3.981268167495728
$ws.Range("E1").Value = 1.784870743751526
